$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.233575666666667
$ws.Range("N2").Value = 12.700727
$ws.Range("O2").Value = 0.4702904532725788
$ws.Range("P2").Value = 0.4702904532725788
$ws.Range("Q2").Value = 114.4712530836954
$ws.Range("R2").Value = 1030.241277753259
$ws.Range("S2").Value = 0.03337269574036905
$ws.Range("T2").Value = 0.03337269574036906

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.605649999999999
$ws.Range("N3").Value = 10.81695
$ws.Range("O3").Value = 0.4005367817548413
$ws.Range("P3").Value = 0.4005367817548413
$ws.Range("Q3").Value = 97.49283021701665
$ws.Range("R3").Value = 877.4354719531498
$ws.Range("S3").Value = 0.02842284391978388
$ws.Range("T3").Value = 0.02842284391978388

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.162819
$ws.Range("N4").Value = 3.488456999999999
$ws.Range("O4").Value = 0.1291727649725799
$ws.Range("P4").Value = 0.1291727649725799
$ws.Range("Q4").Value = 31.44135324840766
$ws.Range("R4").Value = 282.9721792356689
$ws.Range("S4").Value = 0.009166342530184342
$ws.Range("T4").Value = 0.009166342530184344

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.233575666666667
$ws.Range("N5").Value = 12.700727
$ws.Range("O5").Value = 0.4702904532725788
$ws.Range("P5").Value = 0.4702904532725788
$ws.Range("Q5").Value = 1462.982260067644
$ws.Range("R5").Value = 13166.8403406088
$ws.Range("S5").Value = 0.4265146097692985
$ws.Range("T5").Value = 0.4265146097692986

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.605649999999999
$ws.Range("N6").Value = 10.81695
$ws.Range("O6").Value = 0.4005367817548413
$ws.Range("P6").Value = 0.4005367817548413
$ws.Range("Q6").Value = 1245.99213557135
$ws.Range("R6").Value = 11213.92922014215
$ws.Range("S6").Value = 0.3632537891841949
$ws.Range("T6").Value = 0.3632537891841949

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.162819
$ws.Range("N7").Value = 3.488456999999999
$ws.Range("O7").Value = 0.1291727649725799
$ws.Range("P7").Value = 0.1291727649725799
$ws.Range("Q7").Value = 401.8313838262009
$ws.Range("R7").Value = 3616.482454435808
$ws.Range("S7").Value = 0.1171490321815419
$ws.Range("T7").Value = 0.1171490321815419

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.233575666666667
$ws.Range("N8").Value = 12.700727
$ws.Range("O8").Value = 0.4702904532725788
$ws.Range("P8").Value = 0.4702904532725788
$ws.Range("Q8").Value = 35.68370291989245
$ws.Range("R8").Value = 321.1533262790321
$ws.Range("S8").Value = 0.01040314776291121
$ws.Range("T8").Value = 0.01040314776291122

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.605649999999999
$ws.Range("N9").Value = 10.81695
$ws.Range("O9").Value = 0.4005367817548413
$ws.Range("P9").Value = 0.4005367817548413
$ws.Range("Q9").Value = 30.39108157346666
$ws.Range("R9").Value = 273.5197341612
$ws.Range("S9").Value = 0.008860148650862459
$ws.Range("T9").Value = 0.008860148650862463

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.162819
$ws.Range("N10").Value = 3.488456999999999
$ws.Range("O10").Value = 0.1291727649725799
$ws.Range("P10").Value = 0.1291727649725799
$ws.Range("Q10").Value = 9.801097467634666
$ws.Range("R10").Value = 88.209877208712
$ws.Range("S10").Value = 0.002857390260853726
$ws.Range("T10").Value = 0.002857390260853727

